# Refresh the NSE watch-list table (columns B-F, rows 2-57) with the new
# screener results, and append 5 new rows (58-62) that extend the used
# range from A1:F57 to A1:F62.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows need the same "index" style (bold, centered, thin border) that
# column A already carries for rows 2-57, so clone it from row 2 first.
$ws.Cells.Item(2, 1).Copy() | Out-Null
foreach ($r in 58..62) {
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

# Apply every changed cell (row, column) -> new value.
$ws.Cells.Item(2, 2).Value = "NSE:AAATECH"
$ws.Cells.Item(2, 3).Value = "NSE:ADANIPOWER"
$ws.Cells.Item(2, 4).Value = "NSE:BPCL"
$ws.Cells.Item(2, 5).Value = ""
$ws.Cells.Item(3, 2).Value = "NSE:AEROFLEX"
$ws.Cells.Item(3, 3).Value = "NSE:ADVANIHOTR"
$ws.Cells.Item(3, 4).Value = "NSE:CANBK"
$ws.Cells.Item(3, 6).Value = "NSE:DIVISLAB"
$ws.Cells.Item(4, 2).Value = "NSE:AFFLE"
$ws.Cells.Item(4, 3).Value = "NSE:ORIENTCEM"
$ws.Cells.Item(4, 4).Value = "NSE:LAURUSLABS"
$ws.Cells.Item(4, 6).Value = "NSE:HAVELLS"
$ws.Cells.Item(5, 2).Value = "NSE:ALICON"
$ws.Cells.Item(5, 3).Value = "NSE:POKARNA"
$ws.Cells.Item(5, 4).Value = "NSE:MUTHOOTFIN"
$ws.Cells.Item(5, 6).Value = "NSE:IOC"
$ws.Cells.Item(6, 2).Value = "NSE:APTECHT"
$ws.Cells.Item(6, 6).Value = "NSE:IPCALAB"
$ws.Cells.Item(7, 2).Value = "NSE:ARIHANTCAP"
$ws.Cells.Item(7, 6).Value = "NSE:KOTAKBANK"
$ws.Cells.Item(8, 2).Value = "NSE:ASIANPAINT"
$ws.Cells.Item(8, 6).Value = ""
$ws.Cells.Item(9, 2).Value = "NSE:ASIANTILES"
$ws.Cells.Item(9, 6).Value = ""
$ws.Cells.Item(10, 2).Value = "NSE:BANDHANBNK"
$ws.Cells.Item(11, 2).Value = "NSE:BANSWRAS"
$ws.Cells.Item(12, 2).Value = "NSE:BFINVEST"
$ws.Cells.Item(13, 2).Value = "NSE:BIRLACABLE"
$ws.Cells.Item(14, 2).Value = "NSE:BLISSGVS"
$ws.Cells.Item(15, 2).Value = "NSE:CAPTRUST"
$ws.Cells.Item(16, 2).Value = "NSE:CENTURYTEX"
$ws.Cells.Item(17, 2).Value = "NSE:CHAMBLFERT"
$ws.Cells.Item(18, 2).Value = "NSE:CYIENT"
$ws.Cells.Item(19, 2).Value = "NSE:DANGEE"
$ws.Cells.Item(20, 2).Value = "NSE:DEEPINDS"
$ws.Cells.Item(21, 2).Value = "NSE:DEN"
$ws.Cells.Item(22, 2).Value = "NSE:DIAMINESQ"
$ws.Cells.Item(23, 2).Value = "NSE:DUCON"
$ws.Cells.Item(24, 2).Value = "NSE:EIDPARRY"
$ws.Cells.Item(25, 2).Value = "NSE:FDC"
$ws.Cells.Item(26, 2).Value = "NSE:GAEL"
$ws.Cells.Item(27, 2).Value = "NSE:GANESHBE"
$ws.Cells.Item(28, 2).Value = "NSE:GANGESSECU"
$ws.Cells.Item(29, 2).Value = "NSE:GLS"
$ws.Cells.Item(30, 2).Value = "NSE:GOACARBON"
$ws.Cells.Item(31, 2).Value = "NSE:GUJALKALI"
$ws.Cells.Item(32, 2).Value = "NSE:HAPPSTMNDS"
$ws.Cells.Item(33, 2).Value = "NSE:HFCL"
$ws.Cells.Item(34, 2).Value = "NSE:INDIGOPNTS"
$ws.Cells.Item(35, 2).Value = "NSE:INDOCO"
$ws.Cells.Item(36, 2).Value = "NSE:INDOTHAI"
$ws.Cells.Item(37, 2).Value = "NSE:INDRAMEDCO"
$ws.Cells.Item(38, 2).Value = "NSE:IPCALAB"
$ws.Cells.Item(39, 2).Value = "NSE:ITI"
$ws.Cells.Item(40, 2).Value = "NSE:KOTHARIPET"
$ws.Cells.Item(41, 2).Value = "NSE:LAXMICOT"
$ws.Cells.Item(42, 2).Value = "NSE:LEMONTREE"
$ws.Cells.Item(43, 2).Value = "NSE:LYKALABS"
$ws.Cells.Item(44, 2).Value = "NSE:LYPSAGEMS"
$ws.Cells.Item(45, 2).Value = "NSE:MAHASTEEL"
$ws.Cells.Item(46, 2).Value = "NSE:MANAPPURAM"
$ws.Cells.Item(47, 2).Value = "NSE:MANORAMA"
$ws.Cells.Item(48, 2).Value = "NSE:MASTEK"
$ws.Cells.Item(49, 2).Value = "NSE:MIRCELECTR"
$ws.Cells.Item(50, 2).Value = "NSE:MOQUALITY"
$ws.Cells.Item(51, 2).Value = "NSE:MUNJALSHOW"
$ws.Cells.Item(52, 2).Value = "NSE:MUTHOOTCAP"
$ws.Cells.Item(53, 2).Value = "NSE:NILKAMAL"
$ws.Cells.Item(54, 2).Value = "NSE:ORIENTBELL"
$ws.Cells.Item(55, 2).Value = "NSE:PATINTLOG"
$ws.Cells.Item(56, 2).Value = "NSE:PILANIINVS"
$ws.Cells.Item(57, 2).Value = "NSE:PLAZACABLE"
$ws.Cells.Item(58, 1).Value = 56
$ws.Cells.Item(58, 2).Value = "NSE:PNBGILTS"
$ws.Cells.Item(59, 1).Value = 57
$ws.Cells.Item(59, 2).Value = "NSE:PRECWIRE"
$ws.Cells.Item(60, 1).Value = 58
$ws.Cells.Item(60, 2).Value = "NSE:RBLBANK"
$ws.Cells.Item(61, 1).Value = 59
$ws.Cells.Item(61, 2).Value = "NSE:RELIGARE"
$ws.Cells.Item(62, 1).Value = 60
$ws.Cells.Item(62, 2).Value = "NSE:RUCHINFRA"

Write-Host "Updated watch-list: $($ws.UsedRange.Address())"
